$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 295, shifting existing rows 295-344 down to 296-345.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with its data (columns unaffected by the
# insert-copy keep the values that were already carried down from the old
# row 295, but we set every column explicitly to be safe).
$ws.Range("A295").Value2 = 10
$ws.Range("B295").Value2 = "Vega Modelo de Temuco"
$ws.Range("C295").Value2 = "La Araucanía"
$ws.Range("D295").Value2 = 44776
$ws.Range("E295").Value2 = 9
$ws.Range("F295").Value2 = 100112044
$ws.Range("G295").Value2 = "Perejil"
$ws.Range("H295").Value2 = "Sin especificar"
$ws.Range("I295").Value2 = "Primera"
$ws.Range("J295").Value2 = 30
$ws.Range("K295").Value2 = 4000
$ws.Range("L295").Value2 = 4000
$ws.Range("M295").Value2 = 4000
$ws.Range("N295").Value2 = "$/docena de atados (3 kilos)"
$ws.Range("O295").Value2 = "Región Metropolitana"
$ws.Range("P295").Value2 = 1333
$ws.Range("Q295").Value2 = 3
$ws.Range("R295").Value2 = "Hortaliza"
